$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one weekly price record per row (rows 2-113). A new
# weekly record is being added for "Berenjena" at Terminal La Palmera de
# La Serena. It is inserted as the new row 45, pushing every existing
# row from 45 downward by one (45->46, ..., 113->114), which is exactly
# what Excel's row Insert does (including carrying the date number
# format down from the row above into the freshly inserted row).
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly record.
$ws.Cells.Item(45, 1).Value = 8
$ws.Cells.Item(45, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(45, 3).Value = "Coquimbo"
$ws.Cells.Item(45, 4).Value = 44580
$ws.Cells.Item(45, 5).Value = 4
$ws.Cells.Item(45, 6).Value = 100112001
$ws.Cells.Item(45, 7).Value = "Berenjena"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 560
$ws.Cells.Item(45, 11).Value = 8000
$ws.Cells.Item(45, 12).Value = 9000
$ws.Cells.Item(45, 13).Value = 8500
$ws.Cells.Item(45, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(45, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(45, 16).Value = 170
$ws.Cells.Item(45, 17).Value = 50
$ws.Cells.Item(45, 18).Value = "Hortaliza"
